$d = $word.ActiveDocument
$d.Content.Find.Execute("500 equipos de c" + [char]0x00F3 + "mputo por a" + [char]0x00F1 + "o", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 o m" + [char]0x00E1 + "s equipos de c" + [char]0x00F3 + "mputo", 2)
